$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data points in A2:B61 with the new fixed datapoints
$ws.Range("A2").Value = 0.264148405557659
$ws.Range("B2").Value = 0.359414615854502
$ws.Range("A3").Value = 0.271378936151633
$ws.Range("B3").Value = 0.411166529814445
$ws.Range("A4").Value = -0.0655475410000367
$ws.Range("B4").Value = -0.718873814647728
$ws.Range("A5").Value = 0.346458836375797
$ws.Range("B5").Value = 0.431440659983075
$ws.Range("A6").Value = -0.511270562970277
$ws.Range("B6").Value = 0.864376082472041
$ws.Range("A7").Value = 0.311461043984769
$ws.Range("B7").Value = 0.469942300318917
$ws.Range("A8").Value = 0.371938670603871
$ws.Range("B8").Value = 0.369919445979
$ws.Range("A9").Value = 0.312255675652832
$ws.Range("B9").Value = 0.417227361818888
$ws.Range("A10").Value = 0.388181492458542
$ws.Range("B10").Value = 0.529192791874564
$ws.Range("A11").Value = 0.402841643213165
$ws.Range("B11").Value = 0.562090110308997
$ws.Range("A12").Value = 0.409401198293665
$ws.Range("B12").Value = 0.510884539494065
$ws.Range("A13").Value = -0.227913354842788
$ws.Range("B13").Value = 0.860703964612247
$ws.Range("A14").Value = 0.458728217777101
$ws.Range("B14").Value = 0.493670887818288
$ws.Range("A15").Value = 0.10716266871802
$ws.Range("B15").Value = -0.536406312382749
$ws.Range("A16").Value = 0.110610397063658
$ws.Range("B16").Value = -0.475687308016153
$ws.Range("A17").Value = -0.0513563421634618
$ws.Range("B17").Value = 0.846768230829464
$ws.Range("A18").Value = -0.0510442448136673
$ws.Range("B18").Value = 0.876724410767892
$ws.Range("A19").Value = 0.528438207333529
$ws.Range("B19").Value = 0.696790516634549
$ws.Range("A20").Value = 0.295461245528316
$ws.Range("B20").Value = -0.47462798371213
$ws.Range("A21").Value = 0.0209461948990525
$ws.Range("B21").Value = 0.936641695614545
$ws.Range("A22").Value = 0.526950848876572
$ws.Range("B22").Value = 0.494141354139011
$ws.Range("A23").Value = 0.128264128388803
$ws.Range("B23").Value = 0.94446980731873
$ws.Range("A24").Value = 0.584106471697224
$ws.Range("B24").Value = 0.635996373148642
$ws.Range("A25").Value = 0.247443636634555
$ws.Range("B25").Value = -0.36867334299421
$ws.Range("A26").Value = 0.319628885653612
$ws.Range("B26").Value = -0.401886369510798
$ws.Range("A27").Value = 0.385444447974389
$ws.Range("B27").Value = -0.43673526980614
$ws.Range("A28").Value = 0.191497798049511
$ws.Range("B28").Value = 0.835877277400865
$ws.Range("A29").Value = 0.376998998884641
$ws.Range("B29").Value = -0.340993601927503
$ws.Range("A30").Value = 0.356505158765045
$ws.Range("B30").Value = -0.397708641622457
$ws.Range("A31").Value = 0.649649356498543
$ws.Range("B31").Value = 0.630080078582772
$ws.Range("A32").Value = 0.179240247036004
$ws.Range("B32").Value = 0.9425801178901
$ws.Range("A33").Value = 0.590806944667474
$ws.Range("B33").Value = 0.679208941314191
$ws.Range("A34").Value = 0.44552674513274
$ws.Range("B34").Value = -0.267889585595264
$ws.Range("A35").Value = 0.275687287834262
$ws.Range("B35").Value = 0.956766590321791
$ws.Range("A36").Value = 0.732488626812433
$ws.Range("B36").Value = 0.665637964418919
$ws.Range("A37").Value = 0.536094682493435
$ws.Range("B37").Value = -0.211934077033771
$ws.Range("A38").Value = 0.436055082593986
$ws.Range("B38").Value = -0.114185160947856
$ws.Range("A39").Value = 0.548504079624116
$ws.Range("B39").Value = 0.936659462963067
$ws.Range("A40").Value = 0.615326467323541
$ws.Range("B40").Value = -0.154439167892724
$ws.Range("A41").Value = 0.74791794179379
$ws.Range("B41").Value = 0.718970319128864
$ws.Range("A42").Value = 0.546957900737106
$ws.Range("B42").Value = -0.0876125079265752
$ws.Range("A43").Value = 0.48339134177481
$ws.Range("B43").Value = 0.982241450194734
$ws.Range("A44").Value = 0.541895005282746
$ws.Range("B44").Value = -0.0419937146151521
$ws.Range("A45").Value = 0.705566321931483
$ws.Range("B45").Value = 0.735760317406035
$ws.Range("A46").Value = 0.607693707918294
$ws.Range("B46").Value = 1.00440627182139
$ws.Range("A47").Value = 0.546998469756011
$ws.Range("B47").Value = 0.905505912143208
$ws.Range("A48").Value = 0.594617836137195
$ws.Range("B48").Value = 1.01566460038183
$ws.Range("A49").Value = 0.622873532515764
$ws.Range("B49").Value = 0.961216491856635
$ws.Range("A50").Value = 0.81814569543908
$ws.Range("B50").Value = 0.847180492663078
$ws.Range("A51").Value = 0.638391684601757
$ws.Range("B51").Value = 1.00723038864112
$ws.Range("A52").Value = 0.637007476908972
$ws.Range("B52").Value = 1.00234791691739
$ws.Range("A53").Value = 0.52509110614333
$ws.Range("B53").Value = 0.0474836140996691
$ws.Range("A54").Value = 0.512169041194409
$ws.Range("B54").Value = 0.127687162917494
$ws.Range("A55").Value = 0.598538837519143
$ws.Range("B55").Value = 1.02209103726113
$ws.Range("A56").Value = 0.536885816867133
$ws.Range("B56").Value = 0.136950742577372
$ws.Range("A57").Value = 0.599043271394805
$ws.Range("B57").Value = 1.00106266493543
$ws.Range("A58").Value = 0.801738483596803
$ws.Range("B58").Value = 0.792661554506447
$ws.Range("A59").Value = 0.44797841800089
$ws.Range("B59").Value = 0.108373932416459
$ws.Range("A60").Value = 0.60013527650117
$ws.Range("B60").Value = 1.0359568929509
$ws.Range("A61").Value = 0.53113098160141
$ws.Range("B61").Value = 0.203593108204057

# Update the sheet view selection to match the saved state
$ws.Range("A2:B61").Select()
